$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Андрей Фокин"
$ws.Range("B3").Value = "2023-07-17 12:19:35"
